# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and fix the ordering of Toncoin / WrappedBTC (rows 14 and 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.14"
$ws.Range("E14").Value = "  -3.08%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "96.204.46"
$ws.Range("E15").Value = "  -1.11%  "

$ws.Range("D2").Value = "96.520.76"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "3.326.27"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.17"
$ws.Range("E5").Value = "  -2.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "652.29"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.36"
$ws.Range("E7").Value = "  -7.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.415"
$ws.Range("E8").Value = "  -3.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.983"
$ws.Range("E10").Value = "  -6.95%  "
$ws.Range("D11").Value = "3.323.34"
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.205"
$ws.Range("E12").Value = "  -3.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.63"
$ws.Range("E13").Value = "  -5.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000250"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").Value = "3.939.13"
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.42"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "3.323.82"
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.73"
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.477"
$ws.Range("E21").Value = "  -6.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "499.56"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.39"
$ws.Range("E23").Value = "  -3.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.33"
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000197"
$ws.Range("E25").Value = "  -4.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.43"
$ws.Range("E26").Value = "  +4.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "94.24"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.90"
$ws.Range("E28").Value = "  -5.60%  "
$ws.Range("D29").Value = "3.494.17"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.141"
$ws.Range("E31").Value = "  -6.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.74"
$ws.Range("E32").Value = "  -5.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.185"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.46"
$ws.Range("E34").Value = "  +13.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.540"
$ws.Range("E36").Value = "  -4.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.89"
$ws.Range("E37").Value = "  -6.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.50"
$ws.Range("E38").Value = "  +8.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.48"
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.150"
$ws.Range("E41").Value = "  -4.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "503.44"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.47"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.66"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.822"
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0404"
$ws.Range("E46").Value = "  -6.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.29"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.40"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.60"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.66"
$ws.Range("E50").Value = "  +3.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.13"
$ws.Range("E51").Value = "  -5.55%  "
